$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with the new rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.33 = 8771.99 pesos`n✅ 8771.99 pesos = 2.31 = 960.3 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the transfi rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 429.777
$ws2.Range("O10").Value = 3770
$ws2.Range("N12").Value = 3800
$ws2.Range("O12").Value = 416.001
